# Automatic update of files.
# Update B column (Taxonsorteringsordning) for rows 7-10 from 96720 to 96735,
# and swap the A/Q/R/Z/AB values between rows 8 and 9.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7
$ws.Range("B7").Value = 96735

# Row 8 (take over former row 9's A/Q/R/Z/AB values)
$ws.Range("A8").Value = 112092130
$ws.Range("B8").Value = 96735
$ws.Range("Q8").Value = 584352
$ws.Range("R8").Value = 7048232
$ws.Range("Z8").Value = "17:22"
$ws.Range("AB8").Value = "17:22"

# Row 9 (take over former row 8's A/Q/R/Z/AB values)
$ws.Range("A9").Value = 112092066
$ws.Range("B9").Value = 96735
$ws.Range("Q9").Value = 584346
$ws.Range("R9").Value = 7048207
$ws.Range("Z9").Value = "17:18"
$ws.Range("AB9").Value = "17:18"

# Row 10
$ws.Range("B10").Value = 96735
